# Generate Report for Handoff
# Updates the localization-status report to reflect that the first file
# (332823b7-1dea-41b0-9582-f288fc2cfbd8.md) has moved from "In Translation"
# to "Ready for handoff", the priority for all rows changed from "ht" to
# "mt", and the handoff timestamps were refreshed to the new generation
# time.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Overview sheet ---
# Row 2 corresponds to 332823b7-1dea-41b0-9582-f288fc2cfbd8.md
$ws1.Range("E2").Value = "Ready for handoff"
$ws1.Range("F2").Value = "Ready for handoff"

# The "Latest HO Xliff Generate Date" column is refreshed for every row
$ws1.Range("G2").Value = "2017-02-09 17:12:46"
$ws1.Range("G3").Value = "2017-02-09 17:12:46"
$ws1.Range("G4").Value = "2017-02-09 17:12:46"
$ws1.Range("G5").Value = "2017-02-09 17:12:46"

# --- zh-cn sheet ---
$ws2.Range("C2").Value = "Ready for handoff"

$ws2.Range("E2").Value = "mt"
$ws2.Range("E3").Value = "mt"
$ws2.Range("E4").Value = "mt"
$ws2.Range("E5").Value = "mt"

$ws2.Range("H2").Value = "2017-02-09 17:12:29"
$ws2.Range("H3").Value = "2017-02-09 17:12:29"
$ws2.Range("H4").Value = "2017-02-09 17:12:29"
$ws2.Range("H5").Value = "2017-02-09 17:12:29"

# --- de-de sheet ---
$ws3.Range("C2").Value = "Ready for handoff"

$ws3.Range("E2").Value = "mt"
$ws3.Range("E3").Value = "mt"
$ws3.Range("E4").Value = "mt"
$ws3.Range("E5").Value = "mt"

$ws3.Range("H2").Value = "2017-02-09 17:12:46"
$ws3.Range("H3").Value = "2017-02-09 17:12:46"
$ws3.Range("H4").Value = "2017-02-09 17:12:46"
$ws3.Range("H5").Value = "2017-02-09 17:12:46"
